# The commit swaps the deck's colour theme from the custom "Integral" /
# "Red Violet" palette over to the stock PowerPoint "Office Theme" palette
# (the 12 theme colours: dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).
#
# PowerPoint's object model exposes those 12 slots through
# SlideMaster.Theme.ThemeColorScheme (1-based, in DrawingML schema order:
# dk1, lt1, dk2, lt2, accent1..accent6, hlink, folHlink). Each entry's
# .RGB is a normal OLE COLORREF (0x00BBGGRR), so we convert each target
# hex colour before assigning it.

$p = $ppt.ActivePresentation

# Target "Office Theme" colours, in clrScheme schema order.
$officeHex = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

function HexToOleRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# The presentation's single Design / SlideMaster carries the theme that
# actually renders the slides - update its 12 theme colours to the Office
# Theme palette.
$design = $p.Designs.Item(1)
$themeColors = $design.SlideMaster.Theme.ThemeColorScheme

for ($i = 1; $i -le $officeHex.Count; $i++) {
    $themeColors.Item($i).RGB = HexToOleRgb($officeHex[$i - 1])
}
